# Fix equirectangular projection string
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the row in column A that holds "equirec" and fix the proj4 string
# next to it (column B) to the correct equirectangular proj4 definition.
$used = $ws.UsedRange
$rows = $used.Rows.Count

for ($r = 1; $r -le $rows; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    if ($cellA.Value() -eq "equirec") {
        # Leading apostrophe tells Excel to store this as text (it starts
        # with "+", which would otherwise be parsed like a formula) and
        # keeps the same quote-prefix cell style the rest of the column uses.
        $ws.Cells.Item($r, 2).Value = "'+proj=eqc"
        break
    }
}

# Update the active selection to match the saved state (B4)
$ws.Range("B4").Select()
